$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1617
$ws.Range("I9").Value = 1425.5
$ws.Range("J9").Value = 2000
$ws.Range("K9").Value = 1425.5
$ws.Range("L9").Value = 2000
$ws.Range("M9").Value = -1256.5
$ws.Range("N9").Value = -2338

$ws.Range("H20").Value = 3250
$ws.Range("I20").Value = 3000
$ws.Range("J20").Value = 3500
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 3500
$ws.Range("M20").Value = -2770

$ws.Range("H33").Value = 454.8421
$ws.Range("I33").Value = 467.77777
$ws.Range("J33").Value = 222
$ws.Range("K33").Value = 467.77777
$ws.Range("L33").Value = 222
$ws.Range("M33").Value = -238.77777

$ws.Range("H35").Value = 3250
$ws.Range("I35").Value = 3000
$ws.Range("J35").Value = 3500
$ws.Range("K35").Value = 3000
$ws.Range("L35").Value = 3500
$ws.Range("M35").Value = -2621

$ws.Range("H111").Value = 1217.6666
$ws.Range("I111").Value = 1217.6666
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3652.9998
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -585.9998000000001

$ws.Range("H113").Value = 4877.5
$ws.Range("I113").Value = 4877.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4877.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1623.5

$ws.Range("H137").Value = 2207.6775
$ws.Range("I137").Value = 2119.652
$ws.Range("J137").Value = 2460.75
$ws.Range("K137").Value = 6358.956
$ws.Range("L137").Value = 7382.25
$ws.Range("M137").Value = -3808.956
$ws.Range("N137").Value = -12482.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 23957.875
$ws.Range("I63").Value = 30415.75
$ws.Range("J63").Value = 17500
$ws.Range("K63").Value = 30415.75
$ws.Range("L63").Value = 17500
$ws.Range("M63").Value = -29729.75
$ws.Range("N63").Value = -18872

$ws.Range("H66").Value = 23957.875
$ws.Range("I66").Value = 30415.75
$ws.Range("J66").Value = 17500
$ws.Range("K66").Value = 152078.75
$ws.Range("L66").Value = 87500
$ws.Range("M66").Value = -148646.75
$ws.Range("N66").Value = -94364

$ws.Range("H122").Value = 1936.84
$ws.Range("I122").Value = 1936.84
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5810.52
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3360.52
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 6276.8423
$ws.Range("I132").Value = 6309.4116
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 18928.2348
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -16398.2348

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1970.6666
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 2456
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 2456
$ws.Range("M7").Value = -887
$ws.Range("N7").Value = -2682

$ws.Range("H134").Value = 5747.9375
$ws.Range("I134").Value = 6204.5
$ws.Range("J134").Value = 4378.25
$ws.Range("K134").Value = 18613.5
$ws.Range("L134").Value = 13134.75
$ws.Range("M134").Value = -16078.5
$ws.Range("N134").Value = -18204.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 2996
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2996
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2996
$ws.Range("N13").Value = -3274

$ws.Range("H22").Value = 5000484
$ws.Range("I22").Value = 554.2
$ws.Range("J22").Value = 13333701
$ws.Range("K22").Value = 554.2
$ws.Range("L22").Value = 13333701
$ws.Range("M22").Value = -204.2
$ws.Range("N22").Value = -13334401

$ws.Range("H31").Value = 2155.5
$ws.Range("I31").Value = 2155.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2155.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1860.5

$ws.Range("H34").Value = 2155.5
$ws.Range("I34").Value = 2155.5
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2155.5
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1953.5

$ws.Range("H122").Value = 3747.4
$ws.Range("I122").Value = 4549.625
$ws.Range("J122").Value = 2830.5715
$ws.Range("K122").Value = 13648.875
$ws.Range("L122").Value = 8491.7145
$ws.Range("M122").Value = -11198.875
$ws.Range("N122").Value = -13391.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 356.91666
$ws.Range("I12").Value = 389.57144
$ws.Range("J12").Value = 311.2
$ws.Range("K12").Value = 1168.71432
$ws.Range("L12").Value = 933.5999999999999
$ws.Range("M12").Value = -995.71432
$ws.Range("N12").Value = -1279.6

$ws.Range("H18").Value = 368
$ws.Range("I18").Value = 313.42856
$ws.Range("J18").Value = 750
$ws.Range("K18").Value = 940.28568
$ws.Range("L18").Value = 2250
$ws.Range("M18").Value = -771.28568

$ws.Range("H107").Value = 818.25
$ws.Range("I107").Value = 802.25
$ws.Range("J107").Value = 834.25
$ws.Range("K107").Value = 2406.75
$ws.Range("L107").Value = 2502.75
$ws.Range("M107").Value = -486.75
$ws.Range("N107").Value = -6342.75

$ws.Range("H122").Value = 501.625
$ws.Range("I122").Value = 502
$ws.Range("J122").Value = 499
$ws.Range("K122").Value = 4518
$ws.Range("L122").Value = 4491
$ws.Range("M122").Value = -2068
$ws.Range("N122").Value = -9391

$ws.Range("H132").Value = 1222.7142
$ws.Range("I132").Value = 1155.6666
$ws.Range("J132").Value = 1625
$ws.Range("K132").Value = 10400.9994
$ws.Range("L132").Value = 14625
$ws.Range("M132").Value = -7870.999400000001

$ws.Range("H140").Value = 1002702.7
$ws.Range("I140").Value = 1113225.2
$ws.Range("J140").Value = 8000
$ws.Range("K140").Value = 3339675.6
$ws.Range("L140").Value = 24000
$ws.Range("M140").Value = -3334495.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 855.4286
$ws.Range("I31").Value = 623.25
$ws.Range("J31").Value = 1165
$ws.Range("K31").Value = 623.25
$ws.Range("L31").Value = 1165
$ws.Range("M31").Value = -331.25

$ws.Range("H37").Value = 855.4286
$ws.Range("I37").Value = 623.25
$ws.Range("J37").Value = 1165
$ws.Range("K37").Value = 623.25
$ws.Range("L37").Value = 1165
$ws.Range("M37").Value = -346.25

$ws.Range("H70").Value = 18910.166
$ws.Range("I70").Value = 24037.46
$ws.Range("J70").Value = 5579.2
$ws.Range("K70").Value = 24037.46
$ws.Range("L70").Value = 5579.2
$ws.Range("M70").Value = -23767.46
$ws.Range("N70").Value = -6119.2

$ws.Range("H73").Value = 18910.166
$ws.Range("I73").Value = 24037.46
$ws.Range("J73").Value = 5579.2
$ws.Range("K73").Value = 24037.46
$ws.Range("L73").Value = 5579.2
$ws.Range("M73").Value = -23101.46
$ws.Range("N73").Value = -7451.2

$ws.Range("H80").Value = 2767
$ws.Range("I80").Value = 2083.7778
$ws.Range("J80").Value = 4816.6665
$ws.Range("K80").Value = 2083.7778
$ws.Range("L80").Value = 4816.6665
$ws.Range("M80").Value = -1085.7778
$ws.Range("N80").Value = -6812.6665

$ws.Range("H83").Value = 2767
$ws.Range("I83").Value = 2083.7778
$ws.Range("J83").Value = 4816.6665
$ws.Range("K83").Value = 10418.889
$ws.Range("L83").Value = 24083.3325
$ws.Range("M83").Value = -5426.888999999999
$ws.Range("N83").Value = -34067.3325

$ws.Range("H102").Value = 2133.3333
$ws.Range("I102").Value = 2133.3333
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2133.3333
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -511.3332999999998
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 4419
$ws.Range("I122").Value = 3023.75
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 9071.25
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -6621.25
$ws.Range("N122").Value = -34900

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 19000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 19000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 19000
$ws.Range("N2").Value = -19224

$ws.Range("H9").Value = 687
$ws.Range("I9").Value = 687
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 687
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -463

$ws.Range("H22").Value = 1024.875
$ws.Range("I22").Value = 699.8333
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 699.8333
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -404.8333
$ws.Range("N22").Value = -2590

$ws.Range("H27").Value = 1024.875
$ws.Range("I27").Value = 699.8333
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 699.8333
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -592.8333
$ws.Range("N27").Value = -2214

$ws.Range("H35").Value = 1065.25
$ws.Range("I35").Value = 1065.25
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1065.25
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -729.25

$ws.Range("H40").Value = 3476.5
$ws.Range("I40").Value = 3476.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3476.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3340.5

$ws.Range("H130").Value = 66665
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 66665
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 66665
$ws.Range("N130").Value = -76705

$ws.Range("H136").Value = 29413590
$ws.Range("I136").Value = 1595.2142
$ws.Range("J136").Value = 166669570
$ws.Range("K136").Value = 4785.642599999999
$ws.Range("L136").Value = 500008710
$ws.Range("M136").Value = -2235.642599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4181.615
$ws.Range("I81").Value = 2133.818
$ws.Range("J81").Value = 5683.3335
$ws.Range("K81").Value = 4267.636
$ws.Range("L81").Value = 11366.667
$ws.Range("M81").Value = -3206.636
$ws.Range("N81").Value = -13488.667

$ws.Range("H84").Value = 4181.615
$ws.Range("I84").Value = 2133.818
$ws.Range("J84").Value = 5683.3335
$ws.Range("K84").Value = 21338.18
$ws.Range("L84").Value = 56833.335
$ws.Range("M84").Value = -16034.18
$ws.Range("N84").Value = -67441.33499999999

$ws.Range("H122").Value = 3994.7334
$ws.Range("I122").Value = 3786.3845
$ws.Range("J122").Value = 5349
$ws.Range("K122").Value = 11359.1535
$ws.Range("L122").Value = 16047
$ws.Range("M122").Value = -8909.1535
$ws.Range("N122").Value = -20947
